$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Relative - ..." bullet: merge the two runs back into one run and drop
#    the _GoBack bookmark that currently sits at the end of that paragraph
#    (it is relocated to the "Auto inheritance" bullet below).
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$rRelative = $d.Content
$rRelative.Find.Execute( `
    "Relative – The element is moved in relation to where it would have been in the normal flow.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Relative – The element is moved in relation to where it would have been in the normal flow.", `
    2)

# ---------------------------------------------------------------------------
# 2) "Auto inheritance e.g. font-family or color" bullet: the sentence
#    becomes "Auto inheritance usually those related to text e.g.
#    font-family or color", typed in the middle of the original run, which
#    is exactly why the run splits in three and _GoBack now marks the
#    insertion point just before "e.g.".
# ---------------------------------------------------------------------------
$rAuto = $d.Content
$rAuto.Find.Execute( `
    "Auto inheritance e.g. font-family or color", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Auto inheritance usually those related to text e.g. font-family or color", `
    2)

# Split off "Auto inheritance " from the rest with a throwaway bookmark so
# the subsequently-typed phrase lands in its own run instead of being
# absorbed back into the preceding run.
$rAfterAuto = $d.Content
$rAfterAuto.Find.Execute("Auto inheritance ")
$rAfterAuto.Collapse(0)
$d.Bookmarks.Add("zzTmpSplit", $rAfterAuto)

# Place _GoBack right after "usually those related to text " (i.e. right
# before "e.g. font-family or color"), which both marks the real insertion
# point and forces the run split required by the target markup.
$rAfterUsually = $d.Content
$rAfterUsually.Find.Execute("usually those related to text ")
$rAfterUsually.Collapse(0)
$d.Bookmarks.Add("_GoBack", $rAfterUsually)

# Remove the scaffolding bookmark - the run boundary it created stays in
# place even after the bookmark itself is gone.
$d.Bookmarks("zzTmpSplit").Delete()

Write-Output "content edits applied"
